$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "D3"   = -8.718699999999995
    "D14"  = -8.171600000000002
    "D16"  = -8.265499999999999
    "D21"  = -8.416499999999997
    "D23"  = -7.943999999999999
    "D25"  = -7.641599999999999
    "D26"  = -7.361099999999998
    "D29"  = -7.299600000000003
    "D40"  = -8.83979999999999
    "D53"  = -8.000499999999999
    "D57"  = -8.365999999999994
    "D59"  = -8.278699999999994
    "D65"  = -7.804400000000007
    "D69"  = -7.555799999999996
    "D79"  = -7.248400000000004
    "D83"  = -9.064199999999998
    "D91"  = -7.9793
    "D93"  = -6.911999999999991
    "D100" = -7.8552
    "D103" = -7.6767
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
